$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (candidate id 1)
$ws.Range("C2").Value = "The candidate lacks the required skills of MongoDB and NodeJS which are essential for the role. Their past projects show strong skills in Data Analysis and Computer Vision, but these do not directly align with the job requirements."
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "The candidate shows good adaptability by expressing willingness to work in different countries, but may need to work on improving skills independently. The candidate's interest in learning new skills and contributing to cutting-edge technologies aligns with the company's expectations."

# Row 3 (candidate id 2)
$ws.Range("C3").Value = "The candidate lacks the required skills in MongoDB, which is essential for the job. Additionally, the candidate has knowledge of Dart and Firebase, which are not mentioned in the job description."
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "The applicant has a strong motivation to work in Japan, demonstrated by the interest in AI, the desire to be part of Japanese work culture, and the willingness to learn a new language. The applicant also shows adaptability and a collaborative approach, which align with the company's values. However, the preference not to work alone and the need to improve skills may impact the overall personality score."

# Row 4 (candidate id 3)
$ws.Range("C4").Value = "The applicant lacks the required skills of MongoDB and NodeJS, which are crucial for the MERN stack development role. The candidate excels in Python, PyTorch, Tensorflow, and Django Rest Framework, which are not directly relevant to the job description."
$ws.Range("E4").Value = "The applicant demonstrates enthusiasm for AI/ML and shows willingness to adapt to new cultures. They have good teamwork skills and a problem-solving attitude. However, they may need to work on independence and time management."
